# Edit script: 
# 1. Delete the paragraph containing "Cherokee" entirely.
# 2. Remove the existing "_GoBack" bookmark (located before "Changes in class Position").
# 3. Add a new "_GoBack" bookmark at the start of the "Team "Chlorine"" paragraph.

$d = $word.ActiveDocument

# --- Step 1: delete the existing _GoBack bookmark ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: delete the paragraph that contains "Cherokee" ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Cherokee") {
        $p.Range.Delete()
        break
    }
}

# --- Step 3: add the "_GoBack" bookmark at the start of the "Team "Chlorine"" paragraph ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Team*Chlorine*") {
        $start = $p.Range.Start
        $bmRange = $d.Range($start, $start)
        $d.Bookmarks.Add("_GoBack", $bmRange)
        break
    }
}
